$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "max" column (column C). This shifts the old "prediction"
# column (D) into C, and the old "rejection-f" column (E) into D.
$ws.Columns.Item(3).Delete()

# Column B ("1-c__Elusimicrobia") now holds the raw score values instead of
# the placeholder "1" (there is now a single child, so no separate max needed).
$ws.Range("B2").Value = 1643.945200406717
$ws.Range("B3").Value = 1535.092286422305
$ws.Range("B4").Value = 1723.472007440319
